# Updated protocol description templates
#
# Fills in the "situation names" sheet (column A) with the sequential
# situation numbers 1..66 below the existing header row, then makes
# "situation names" the active sheet/tab with that range selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("situation names")

$count = 66
$values = New-Object 'object[,]' $count,1
for ($i = 1; $i -le $count; $i++) {
    $values[$i - 1, 0] = $i
}
$ws.Range("A2:A67").Value = $values

$ws.Activate()
$ws.Range("A2:A67").Select() | Out-Null
